$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): Fiid | Tmice | Th | NaK_ratio
#   - "Fiid", "Th", "Tmice" already existed as shared strings in the workbook
#   - "NaK_ratio" replaces the old (now-unused) "DD01_0x_0x" id strings
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Fiid"
$ws.Range("B1").Value = "Tmice"
$ws.Range("C1").Value = "Th"
$ws.Range("D1").Value = "NaK_ratio"

# ---------------------------------------------------------------------------
# Data rows (2 .. 23): 22 observations.
#   A = Fiid running id 1..22
#   B = Tmice = -Fiid
#   C = Th    = measured value
#   D = NaK_ratio, constant 1 for every row
# ---------------------------------------------------------------------------
$th = @(100, 109, 110, 80, 120, 90, 85, 105, 125, 140, 195, 165, 140, 98, 160, 85, 142, 180, 75, 68, 91, 75)

for ($i = 0; $i -lt $th.Count; $i++) {
    $row  = $i + 2
    $fiid = $i + 1

    $ws.Cells.Item($row, 1).Value = $fiid
    $ws.Cells.Item($row, 2).Value = (0 - $fiid)
    $ws.Cells.Item($row, 3).Value = $th[$i]
    $ws.Cells.Item($row, 4).Value = 1
}

# ---------------------------------------------------------------------------
# View state: scroll so row 6 is at the top and select C24 (next empty row),
# matching the saved sheetView/selection in the target workbook.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C24").Select()
